function Wrap-Xml($body) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $footer = '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + '<w:body>' + $body + '</w:body>' + $footer
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Block 1: "Target audience" / "Problem Statement" / "This could be used..."
# paragraphs -> add intro paragraph, bold the labels, split runs, add a new
# closing paragraph with the moved _GoBack bookmark.
# ---------------------------------------------------------------------------
$block1 = @'
<w:p>
  <w:r>
    <w:t>In 2018, there was a total of 160,597 casualties of all severities in road traffic crashes</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">. Of these 1,784 were fatal and 25,511 were serious injuries. On average 5 people die every day in Great Britain, the reduction in deaths from 2017 to 2018 was just 1%, so I wanted to investigate what conditions lead to the severity of an accident increasing and whether a model could be created to both predict the likelihood of a serious crash to allow for increased availability of emergency services to hotspot areas and to model how temporary changes (such as roadworks or speed restrictions) can impact this. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Target audience:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Local authorities</w:t>
  </w:r>
  <w:r>
    <w:t>/Emergency services</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> in UK</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Problem Statement:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>predict likelihood of severe accidents on roads on a specific day</w:t>
  </w:r>
  <w:r>
    <w:t>/time</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> given the weather conditions. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">This could be used </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">by local authorities </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">to model the </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">current risk and then the </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">reduced risk if temporary measures are put in place (i.e. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">create a prediction with </w:t>
  </w:r>
  <w:r>
    <w:t>lower speed limits)</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">, therefore allowing </w:t>
  </w:r>
  <w:r>
    <w:t>the most effective measures</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> to be put in place to prevent accidents from occurrin</w:t>
  </w:r>
  <w:r>
    <w:t>g</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Emergency services could also look at the likelihood of severe incidents on different road classes/areas, to allow for advanced planning of resources.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$pStart = $d.Paragraphs(4)
$pEnd = $d.Paragraphs(6)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rng.InsertXML((Wrap-Xml $block1))

Write-Output "Block1 done. Paragraph count: $($d.Paragraphs.Count)"
